$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: 309 -> 346
$ws.Range("B1").Value = 346.0

# A4: "20161111" -> "20161122" (keep as text, not auto-converted to a number)
$ws.Range("A4").Value = "'20161122"

# Append new trailing rows (59-66) found in the updated data set
$newRows = @(
    @{ A = "20200812"; B = 1165.0 },
    @{ A = "20201013"; B = 2534.0 },
    @{ A = "20210114"; B = 1199.0 },
    @{ A = "20210120"; B = 1403.0 },
    @{ A = "20210208"; B = 1066.0 },
    @{ A = "20210302"; B = 1264.0 },
    @{ A = "20210315"; B = 1077.0 },
    @{ A = "20210402"; B = 1269.0 }
)

$startRow = 59
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "'" + $newRows[$i].A
    $ws.Cells.Item($r, 2).Value = $newRows[$i].B
}
